# "yeni dosya sisteminde ilk değişikilik"
# Add a second value "deneme1" next to the existing "deneme" cell (A1),
# in B1, and leave the selection on B2 (matching the saved workbook view).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = "deneme1"
$ws.Range("B2").Select()
